$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append two new data rows (101 and 102) to the YKY.MI price/volume table,
# continuing the existing pattern: date serial in A (same date-time style as
# the rows above), volume=0, O/H/L/C=1, G holds the text "1", H the ticker.

$ws.Cells.Item(101, 1).Value = 45454.2916666667
$ws.Cells.Item(101, 2).Value = 0
$ws.Cells.Item(101, 3).Value = 1
$ws.Cells.Item(101, 4).Value = 1
$ws.Cells.Item(101, 5).Value = 1
$ws.Cells.Item(101, 6).Value = 1
$ws.Cells.Item(101, 8).Value = "YKY.MI"

$ws.Cells.Item(102, 1).Value = 45455.2916666667
$ws.Cells.Item(102, 2).Value = 0
$ws.Cells.Item(102, 3).Value = 1
$ws.Cells.Item(102, 4).Value = 1
$ws.Cells.Item(102, 5).Value = 1
$ws.Cells.Item(102, 6).Value = 1
$ws.Cells.Item(102, 8).Value = "YKY.MI"

# Reuse the existing date/time number format (style) from the row above for
# column A, and the existing text "1" (already stored as a shared string) for
# column G, by copying the formatted/typed source cells down, one row at a
# time (paste does not tile across a multi-cell destination).
$ws.Range("A100").Copy()
$ws.Range("A101").PasteSpecial(-4122)
$ws.Range("A100").Copy()
$ws.Range("A102").PasteSpecial(-4122)

$ws.Range("G100").Copy()
$ws.Range("G101").PasteSpecial(-4163)
$ws.Range("G100").Copy()
$ws.Range("G102").PasteSpecial(-4163)
